# Apply updated crypto price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.964.70"
$ws.Range("E2").Value = "  +0.96%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.950.91"
$ws.Range("E3").Value = "  -0.74%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.13"
$ws.Range("E5").Value = "  -1.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4890"
$ws.Range("E7").Value = "  +1.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2965"
$ws.Range("E8").Value = "  +0.65%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06826"
$ws.Range("E9").Value = "  +0.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.14"
$ws.Range("E10").Value = "  -1.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "107.39"
$ws.Range("E11").Value = "  -3.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.956.34"
$ws.Range("E12").Value = "  -0.51%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07778"
$ws.Range("E13").Value = "  +0.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.451"
$ws.Range("E14").Value = "  -0.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7040"
$ws.Range("E15").Value = "  +1.83%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "282.83"
$ws.Range("E16").Value = "  -3.61%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "31.011.11"
$ws.Range("E17").Value = "  +1.10%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.22"
$ws.Range("E18").Value = "  -0.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007698"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.208.28"
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.478"
$ws.Range("E22").Value = "  -3.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9987"
$ws.Range("E23").Value = "  -0.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.474"
$ws.Range("E24").Value = "  -2.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.831"
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.77"
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.96"
$ws.Range("E27").Value = "  -1.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.196"
$ws.Range("E28").Value = "  -0.98%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1057"
$ws.Range("E29").Value = "  -2.54%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.410"
$ws.Range("E30").Value = "  -2.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.579"
$ws.Range("E31").Value = "  -1.63%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.591"
$ws.Range("E32").Value = "  -1.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.442"
$ws.Range("E33").Value = "  +0.29%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04951"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7657"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.169"
$ws.Range("E36").Value = "  -0.52%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.726"
$ws.Range("E37").Value = "  -0.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02011"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.700"
$ws.Range("E39").Value = "  -0.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.574"
$ws.Range("E40").Value = "  +8.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.141"
$ws.Range("E41").Value = "  +3.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "73.71"
$ws.Range("E42").Value = "  +5.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4480"
$ws.Range("E43").Value = "  +0.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "109.64"
$ws.Range("E44").Value = "  -1.63%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8845"
$ws.Range("E45").Value = "  +0.95%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.112"
$ws.Range("E46").Value = "  +9.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9991"
$ws.Range("E47").Value = "  -0.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "990.61"
$ws.Range("E48").Value = "  +7.83%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1262"
$ws.Range("E49").Value = "  +0.18%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.312"
$ws.Range("E50").Value = "  -0.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "35.70"
$ws.Range("E51").Value = "  -0.45%  "
